{"js": "// The document used a series of \"horizontal rule\" paragraphs (a lone\n// run containing a VML <w:pict><v:rect .../></w:pict> divider, with no\n// visible text) to separate the write-up into sections. Those divider\n// paragraphs are being removed now that the outline's sub-sections are\n// fleshed out with real headings, so headings/body text now sit flush\n// against one another instead of being separated by a rule.\n//\n// A divider paragraph is fully empty as far as `Paragraph.text` is\n// concerned (the picture contributes no text), which is not true for\n// any other paragraph in this document, so that is a safe, precise way\n// to find them without depending on fragile paragraph indices.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst dividers = paragraphs.items.filter((p) => p.text.length === 0);\ndividers.forEach((p) => p.delete());\n\nawait context.sync();\n", "ps1": "# The document used a series of \"horizontal rule\" paragraphs (a lone\n# run containing a VML <w:pict><v:rect .../></w:pict> divider, with no\n# visible text) to separate the write-up into sections. Those divider\n# paragraphs are being removed now that the outline's sub-sections are\n# fleshed out with real headings, so headings/body text now sit flush\n# against one another instead of being separated by a rule.\n#\n# In the Word object model, Range.Text for such a paragraph is just the\n# paragraph mark (length 1, no visible characters) because the picture\n# contributes no text - that is not true for any other paragraph in this\n# document, so trimming and checking for an empty string is a safe,\n# precise way to find them without depending on fragile paragraph\n# indices. Walk backwards so deleting a paragraph never invalidates the\n# index of paragraphs still to be visited.\n\n$d = $word.ActiveDocument\n$paragraphs = $d.Paragraphs\n$count = $paragraphs.Count\n\nfor ($i = $count; $i -ge 1; $i--) {\n    $para = $paragraphs.Item($i)\n    if ($para.Range.Text.Trim().Length -eq 0) {\n        $para.Range.Delete()\n    }\n}\n"}
